# Update dags to reflect new urbanicity
# Replace the single "C:Urbanicity" / urbanscore_cont_clst row with four
# new urbanicity covariate rows (built_population_2014, nightlights_composite,
# all_population_count_2015, travel_times_2015), highlighted with the new
# "Accent5, Darker 25%" fill + explicit-black font used to flag the updated
# block, and move the selection onto the newly inserted B column range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: push everything from the old "I:Pv18s" row (19) down by 3 rows,
# copying row 18's formatting into the newly inserted rows 19-21.
$ws.Rows("19:21").Insert()

# Row 18 keeps its "C:Urbanicity" label (A18) but gets the first of the four
# new dhs covariate names.
$ws.Range("B18").Value = "built_population_2014_cont_scale_clst"

# The three freshly inserted rows also describe "C:Urbanicity", each paired
# with one of the other new covariates.
$ws.Range("A19").Value = "C:Urbanicity"
$ws.Range("B19").Value = " nightlights_composite_cont_scale_clst"

$ws.Range("A20").Value = "C:Urbanicity"
$ws.Range("B20").Value = " all_population_count_2015_cont_scale_clst"

$ws.Range("A21").Value = "C:Urbanicity"
$ws.Range("B21").Value = " travel_times_2015_cont_scale_clst"

# Highlight the whole updated urbanicity block (A18:B21) with the new fill
# (theme Accent5, darker 25% ~= RGB 2E75B6) and explicit black font so it
# stands out as the newly revised section.
$highlightColor = 11957550  # RGB(0x2E,0x75,0xB6) packed as BGR OLE color
18..21 | ForEach-Object {
    $ws.Range("A$_").Font.Color = 0
    $ws.Range("A$_").Interior.Color = $highlightColor
    $ws.Range("B$_").Interior.Color = $highlightColor
}

# Reflect the edit location in the sheet's active selection, matching where
# the author was last working.
$ws.Range("B18:B21").Select() | Out-Null
